$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 70
$ws_ALC.Range("H70").Value = 1986.2142
$ws_ALC.Range("I70").Value = 1086
$ws_ALC.Range("K70").Value = 3258
$ws_ALC.Range("M70").Value = -2988
# Row 73
$ws_ALC.Range("H73").Value = 1986.2142
$ws_ALC.Range("I73").Value = 1086
$ws_ALC.Range("K73").Value = 3258
$ws_ALC.Range("M73").Value = -2322
# Row 98
$ws_ALC.Range("H98").Value = 13139053
$ws_ALC.Range("I98").Value = 24084380
$ws_ALC.Range("J98").Value = 4659.8
$ws_ALC.Range("K98").Value = 24084380
$ws_ALC.Range("L98").Value = 4659.8
$ws_ALC.Range("M98").Value = -24082882
$ws_ALC.Range("N98").Value = -7655.8
# Row 100
$ws_ALC.Range("H100").Value = 14494361
$ws_ALC.Range("J100").Value = 3770
$ws_ALC.Range("L100").Value = 3770
$ws_ALC.Range("N100").Value = -4852
# Row 113
$ws_ALC.Range("H113").Value = 3188.75
$ws_ALC.Range("I113").Value = 3144.2856
$ws_ALC.Range("J113").Value = 3500
$ws_ALC.Range("K113").Value = 3144.2856
$ws_ALC.Range("L113").Value = 3500
$ws_ALC.Range("M113").Value = 109.7143999999998
$ws_ALC.Range("N113").Value = -10008
# Row 122
$ws_ALC.Range("H122").Value = 13139053
$ws_ALC.Range("I122").Value = 24084380
$ws_ALC.Range("J122").Value = 4659.8
$ws_ALC.Range("K122").Value = 72253140
$ws_ALC.Range("L122").Value = 13979.4
$ws_ALC.Range("M122").Value = -72250690
$ws_ALC.Range("N122").Value = -18879.4
# Row 137
$ws_ALC.Range("H137").Value = 1048.8206
$ws_ALC.Range("I137").Value = 940.37933
$ws_ALC.Range("J137").Value = 1363.3
$ws_ALC.Range("K137").Value = 2821.13799
$ws_ALC.Range("L137").Value = 4089.9
$ws_ALC.Range("M137").Value = -271.1379900000002
$ws_ALC.Range("N137").Value = -9189.9
# Row 138
$ws_ALC.Range("H138").Value = 4178.7256
$ws_ALC.Range("I138").Value = 1824.0333
$ws_ALC.Range("J138").Value = 7542.5713
$ws_ALC.Range("K138").Value = 5472.0999
$ws_ALC.Range("L138").Value = 22627.7139
$ws_ALC.Range("M138").Value = -332.0999000000002
$ws_ALC.Range("N138").Value = -32907.7139
# Row 141
$ws_ALC.Range("H141").Value = 3231.5
$ws_ALC.Range("I141").Value = 3371.4614
$ws_ALC.Range("K141").Value = 10114.3842
$ws_ALC.Range("M141").Value = -4934.3842

# ---- Sheet: ARM ----
$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 32
$ws_ARM.Range("H32").Value = 337821.72
$ws_ARM.Range("I32").Value = 2990.5195
$ws_ARM.Range("J32").Value = 3202488.8
$ws_ARM.Range("K32").Value = 2990.5195
$ws_ARM.Range("L32").Value = 3202488.8
$ws_ARM.Range("M32").Value = -2703.5195
$ws_ARM.Range("N32").Value = -3203062.8
# Row 45
$ws_ARM.Range("H45").Value = 2918.1785
$ws_ARM.Range("I45").Value = 2967.8333
$ws_ARM.Range("J45").Value = 2828.8
$ws_ARM.Range("K45").Value = 2967.8333
$ws_ARM.Range("L45").Value = 2828.8
$ws_ARM.Range("M45").Value = -2590.8333
$ws_ARM.Range("N45").Value = -3582.8
# Row 74
$ws_ARM.Range("H74").Value = 882.7037
$ws_ARM.Range("I74").Value = 750
$ws_ARM.Range("J74").Value = 1261.8572
$ws_ARM.Range("K74").Value = 750
$ws_ARM.Range("L74").Value = 1261.8572
$ws_ARM.Range("M74").Value = 124
$ws_ARM.Range("N74").Value = -3009.8572
# Row 77
$ws_ARM.Range("H77").Value = 882.7037
$ws_ARM.Range("I77").Value = 750
$ws_ARM.Range("J77").Value = 1261.8572
$ws_ARM.Range("K77").Value = 3750
$ws_ARM.Range("L77").Value = 6309.286
$ws_ARM.Range("M77").Value = 618
$ws_ARM.Range("N77").Value = -15045.286
# Row 122
$ws_ARM.Range("H122").Value = 13160928
$ws_ARM.Range("I122").Value = 3325.9092
$ws_ARM.Range("J122").Value = 31252632
$ws_ARM.Range("K122").Value = 9977.7276
$ws_ARM.Range("L122").Value = 93757896
$ws_ARM.Range("M122").Value = -7527.7276
$ws_ARM.Range("N122").Value = -93762796
# Row 132
$ws_ARM.Range("H132").Value = 22774428
$ws_ARM.Range("I132").Value = 27779238
$ws_ARM.Range("J132").Value = 252781.5
$ws_ARM.Range("K132").Value = 83337714
$ws_ARM.Range("L132").Value = 758344.5
$ws_ARM.Range("M132").Value = -83335184
$ws_ARM.Range("N132").Value = -763404.5

# ---- Sheet: BSM ----
$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 105
$ws_BSM.Range("H105").Value = 3500.9473
$ws_BSM.Range("I105").Value = 3035.8333
$ws_BSM.Range("J105").Value = 4298.2856
$ws_BSM.Range("K105").Value = 3035.8333
$ws_BSM.Range("L105").Value = 4298.2856
$ws_BSM.Range("M105").Value = -1288.8333
$ws_BSM.Range("N105").Value = -7792.2856

# ---- Sheet: CRP ----
$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 51
$ws_CRP.Range("H51").Value = 13666.333
$ws_CRP.Range("I51").Value = 3000
$ws_CRP.Range("J51").Value = 18999.5
$ws_CRP.Range("K51").Value = 3000
$ws_CRP.Range("L51").Value = 18999.5
$ws_CRP.Range("M51").Value = -2264
$ws_CRP.Range("N51").Value = -20471.5
# Row 59
$ws_CRP.Range("H59").Value = 11084.667
$ws_CRP.Range("J59").Value = 12127
$ws_CRP.Range("L59").Value = 12127
$ws_CRP.Range("N59").Value = -14417
# Row 61
$ws_CRP.Range("H61").Value = 13666.333
$ws_CRP.Range("I61").Value = 3000
$ws_CRP.Range("J61").Value = 18999.5
$ws_CRP.Range("K61").Value = 3000
$ws_CRP.Range("L61").Value = 18999.5
$ws_CRP.Range("M61").Value = -2652
$ws_CRP.Range("N61").Value = -19695.5
# Row 99
$ws_CRP.Range("H99").Value = 1004089.06
$ws_CRP.Range("I99").Value = 1004089.06
$ws_CRP.Range("K99").Value = 1004089.06
$ws_CRP.Range("M99").Value = -1002591.06
# Row 107
$ws_CRP.Range("H107").Value = 579.46155
$ws_CRP.Range("I107").Value = 585.5454999999999
$ws_CRP.Range("J107").Value = 546
$ws_CRP.Range("K107").Value = 585.5454999999999
$ws_CRP.Range("L107").Value = 546
$ws_CRP.Range("M107").Value = 1334.4545
$ws_CRP.Range("N107").Value = -4386
# Row 122
$ws_CRP.Range("H122").Value = 11765698
$ws_CRP.Range("I122").Value = 884
$ws_CRP.Range("J122").Value = 33334524
$ws_CRP.Range("K122").Value = 2652
$ws_CRP.Range("L122").Value = 100003572
$ws_CRP.Range("M122").Value = -202
$ws_CRP.Range("N122").Value = -100008472
# Row 126
$ws_CRP.Range("H126").Value = 1004089.06
$ws_CRP.Range("I126").Value = 1004089.06
$ws_CRP.Range("K126").Value = 3012267.18
$ws_CRP.Range("M126").Value = -3009797.18
# Row 134
$ws_CRP.Range("H134").Value = 1896.8286
$ws_CRP.Range("I134").Value = 1837.7307
$ws_CRP.Range("K134").Value = 5513.1921
$ws_CRP.Range("M134").Value = -2978.1921

# ---- Sheet: CUL ----
$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 80
$ws_CUL.Range("H80").Value = 2500
$ws_CUL.Range("J80").Value = 2500
$ws_CUL.Range("L80").Value = 7500
$ws_CUL.Range("N80").Value = -9372
# Row 83
$ws_CUL.Range("H83").Value = 2500
$ws_CUL.Range("J83").Value = 2500
$ws_CUL.Range("L83").Value = 22500
$ws_CUL.Range("N83").Value = -31860
# Row 113
$ws_CUL.Range("H113").Value = 904
$ws_CUL.Range("I113").Value = 873
$ws_CUL.Range("J113").Value = 907.3905999999999
$ws_CUL.Range("K113").Value = 2619
$ws_CUL.Range("L113").Value = 2722.1718
$ws_CUL.Range("M113").Value = -449
$ws_CUL.Range("N113").Value = -7062.1718
# Row 132
$ws_CUL.Range("H132").Value = 980.5454999999999
$ws_CUL.Range("I132").Value = 714
$ws_CUL.Range("J132").Value = 1104.9333
$ws_CUL.Range("K132").Value = 6426
$ws_CUL.Range("L132").Value = 9944.3997
$ws_CUL.Range("M132").Value = -3896
$ws_CUL.Range("N132").Value = -15004.3997

# ---- Sheet: GSM ----
$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 122
$ws_GSM.Range("H122").Value = 2877.4
$ws_GSM.Range("I122").Value = 2346.75
$ws_GSM.Range("K122").Value = 7040.25
$ws_GSM.Range("M122").Value = -4590.25
# Row 126
$ws_GSM.Range("H126").Value = 15165399
$ws_GSM.Range("I126").Value = 14964
$ws_GSM.Range("J126").Value = 55566556
$ws_GSM.Range("K126").Value = 44892
$ws_GSM.Range("L126").Value = 166699668
$ws_GSM.Range("M126").Value = -42422
$ws_GSM.Range("N126").Value = -166704608

# ---- Sheet: LTW ----
$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 7
$ws_LTW.Range("H7").Value = 3014
$ws_LTW.Range("I7").Value = 2919.8
$ws_LTW.Range("J7").Value = 3249.5
$ws_LTW.Range("K7").Value = 2919.8
$ws_LTW.Range("L7").Value = 3249.5
$ws_LTW.Range("M7").Value = -2807.8
$ws_LTW.Range("N7").Value = -3473.5
# Row 40
$ws_LTW.Range("H40").Value = 11767059
$ws_LTW.Range("I40").Value = 2366.5833
$ws_LTW.Range("J40").Value = 40002320
$ws_LTW.Range("K40").Value = 2366.5833
$ws_LTW.Range("L40").Value = 40002320
$ws_LTW.Range("M40").Value = -2230.5833
$ws_LTW.Range("N40").Value = -40002592
# Row 46
$ws_LTW.Range("H46").Value = 25700.25
$ws_LTW.Range("I46").Value = 933.6667
$ws_LTW.Range("J46").Value = 100000
$ws_LTW.Range("K46").Value = 933.6667
$ws_LTW.Range("L46").Value = 100000
$ws_LTW.Range("M46").Value = -745.6667
$ws_LTW.Range("N46").Value = -100376
# Row 61
$ws_LTW.Range("H61").Value = 2469.4614
$ws_LTW.Range("I61").Value = 1014
$ws_LTW.Range("J61").Value = 4167.5
$ws_LTW.Range("K61").Value = 1014
$ws_LTW.Range("L61").Value = 4167.5
$ws_LTW.Range("M61").Value = -812
$ws_LTW.Range("N61").Value = -4571.5
# Row 100
$ws_LTW.Range("H100").Value = 2266.7454
$ws_LTW.Range("I100").Value = 1133.4166
$ws_LTW.Range("K100").Value = 1133.4166
$ws_LTW.Range("M100").Value = -592.4166
# Row 113
$ws_LTW.Range("H113").Value = 2469.4614
$ws_LTW.Range("I113").Value = 1014
$ws_LTW.Range("J113").Value = 4167.5
$ws_LTW.Range("K113").Value = 1014
$ws_LTW.Range("L113").Value = 4167.5
$ws_LTW.Range("M113").Value = 1156
$ws_LTW.Range("N113").Value = -8507.5
# Row 126
$ws_LTW.Range("H126").Value = 3014
$ws_LTW.Range("I126").Value = 2919.8
$ws_LTW.Range("J126").Value = 3249.5
$ws_LTW.Range("K126").Value = 8759.400000000001
$ws_LTW.Range("L126").Value = 9748.5
$ws_LTW.Range("M126").Value = -6289.400000000001
$ws_LTW.Range("N126").Value = -14688.5
# Row 136
$ws_LTW.Range("H136").Value = 35752744
$ws_LTW.Range("I136").Value = 44736.824
$ws_LTW.Range("J136").Value = 200009580
$ws_LTW.Range("K136").Value = 134210.472
$ws_LTW.Range("L136").Value = 600028740
$ws_LTW.Range("M136").Value = -131660.472
$ws_LTW.Range("N136").Value = -600033840

# ---- Sheet: WVR ----
$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 113
$ws_WVR.Range("H113").Value = 495.77777
$ws_WVR.Range("I113").Value = 495.77777
$ws_WVR.Range("J113").Value = 0
$ws_WVR.Range("K113").Value = 1487.33331
$ws_WVR.Range("L113").Value = 0
$ws_WVR.Range("M113").Value = 682.66669
$ws_WVR.Range("N113").ClearContents()
# Row 126
$ws_WVR.Range("H126").Value = 1488.25
$ws_WVR.Range("I126").Value = 1339.8889
$ws_WVR.Range("J126").Value = 1933.3334
$ws_WVR.Range("K126").Value = 4019.6667
$ws_WVR.Range("L126").Value = 5800.0002
$ws_WVR.Range("M126").Value = -1549.6667
$ws_WVR.Range("N126").Value = -10740.0002
# Row 132
$ws_WVR.Range("H132").Value = 43273120
$ws_WVR.Range("I132").Value = 77587850
$ws_WVR.Range("J132").Value = 6721.913
$ws_WVR.Range("K132").Value = 232763550
$ws_WVR.Range("L132").Value = 20165.739
$ws_WVR.Range("M132").Value = -232761020
$ws_WVR.Range("N132").Value = -25225.739
